$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new data row (row 4) mirroring the existing rows (2 and 3)
$ws.Range("A4").Value = "2019-04-09 16:22:16"
$ws.Range("B4").Value = "NODE1"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = 1
